# Apply the edits described by the diff for sheet1 ("7.3.1.1" table):
#  - Add a new year column S (2022) with data point 13.5
#  - Update a few existing 2019/2020/2021 values in the data row
#  - Update the sheet selection to reflect the newly added column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column S, row 4 (year header) ---
# Copy formatting from R4 (style index 19) onto S4, then set its value.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(4, 19).Value = 2022

# --- Add new column S, row 5 (data value) ---
# Copy formatting from R5 (style index 8) onto S5, then set its value.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(5, 19).Value = 13.5

# Clear the marching-ants clipboard marquee left over from the copies above.
$ws.Application.CutCopyMode = $false

# --- Update existing data values in row 5 ---
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.899999999999999

# --- Update the stored selection to match the new active cell/range ---
[void]$ws.Range("S7:S8").Select()
